# Apply the two changes captured by the commit:
#  1. Slide 5's table switches from the custom "Table_0" style
#     ({F1453A40-4581-40EB-B814-7B99970578AB}) to the built-in
#     "No Style, Table Grid" style ({2A13C0FA-8BB3-4102-9A73-B4DD16422824}).
#  2. The deck's applied theme (color scheme) reverts from the "Integral"
#     palette to the default "Office Theme" palette (dk1/lt1/dk2/lt2/
#     accent1-6/hlink/folHlink) -- i.e. the theme previously used by the
#     notes master becomes the one used by the slide master (and vice
#     versa), which in the saved OOXML shows up as the contents of
#     theme1.xml/theme2.xml being swapped while the files keep their
#     names.

$p = $ppt.ActivePresentation

# --- 1. Fix the table style on the one table in the deck -----------------
$oldStyleId = "{F1453A40-4581-40EB-B814-7B99970578AB}"
$newStyleId = "{2A13C0FA-8BB3-4102-9A73-B4DD16422824}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Restore the default "Office Theme" color scheme ------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (in that order), expressed
# as COM RGB integers (0x00BBGGRR).
$officeThemeRgb = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slideForTheme = $p.Slides.Item(1)
$themeColors = $slideForTheme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRgb[$i - 1]
}
